# Workbook: "Hortaliza, Vega Monumental Concepción - Acelga"
# Insert two new weekly observation rows ("Primera" / "Segunda") at the top
# of the Acelga data block (rows 319-320), pushing the existing rows 319-447
# down to 321-449.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows at row 319 (everything from old row 319 onward
# shifts down by 2; the sheet's dimension grows from R447 to R449).
$ws.Rows.Item(319).Resize(2).Insert()

# ---- New row 319: Calidad "Primera" ----
$ws.Cells.Item(319, 1).Value  = 11
$ws.Cells.Item(319, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(319, 3).Value  = "Bíobío"
$ws.Cells.Item(319, 4).Value  = 45141
$ws.Cells.Item(319, 5).Value  = 8
$ws.Cells.Item(319, 6).Value  = 100112009
$ws.Cells.Item(319, 7).Value  = "Acelga"
$ws.Cells.Item(319, 8).Value  = "Sin especificar"
$ws.Cells.Item(319, 9).Value  = "Primera"
$ws.Cells.Item(319, 10).Value = 200
$ws.Cells.Item(319, 11).Value = 600
$ws.Cells.Item(319, 12).Value = 700
$ws.Cells.Item(319, 13).Value = 650
$ws.Cells.Item(319, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(319, 15).Value = "Región de Ñuble"
$ws.Cells.Item(319, 16).Value = 650
$ws.Cells.Item(319, 17).Value = 1
$ws.Cells.Item(319, 18).Value = "Hortaliza"

# ---- New row 320: Calidad "Segunda" ----
$ws.Cells.Item(320, 1).Value  = 11
$ws.Cells.Item(320, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(320, 3).Value  = "Bíobío"
$ws.Cells.Item(320, 4).Value  = 45141
$ws.Cells.Item(320, 5).Value  = 8
$ws.Cells.Item(320, 6).Value  = 100112009
$ws.Cells.Item(320, 7).Value  = "Acelga"
$ws.Cells.Item(320, 8).Value  = "Sin especificar"
$ws.Cells.Item(320, 9).Value  = "Segunda"
$ws.Cells.Item(320, 10).Value = 100
$ws.Cells.Item(320, 11).Value = 500
$ws.Cells.Item(320, 12).Value = 500
$ws.Cells.Item(320, 13).Value = 500
$ws.Cells.Item(320, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(320, 15).Value = "Región de Ñuble"
$ws.Cells.Item(320, 16).Value = 500
$ws.Cells.Item(320, 17).Value = 1
$ws.Cells.Item(320, 18).Value = "Hortaliza"
